# Batch-run related changes to the TestPlan sheet:
#  - Drop the old OpCo/Execute columns, keep S.NO/Test_Case as the first two
#    columns followed by the Bank Country/Bank key/Currency Code/Partner
#    Bank Type columns.
#  - Replace the single sample row with five concrete test-script rows that
#    all share the same Panama bank reference data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestPlan")

# Header row
$ws.Cells.Item(1, 1).Value = "S.NO"
$ws.Cells.Item(1, 2).Value = "Test_Case"
$ws.Cells.Item(1, 3).Value = "Bank Country"
$ws.Cells.Item(1, 4).Value = "Bank key"
$ws.Cells.Item(1, 5).Value = "Currency Code"
$ws.Cells.Item(1, 6).Value = "Partner Bank Type"

# Clear the old trailing columns (G:P) that used to hold extra header/blank cells
$ws.Range("G1:P1").Clear()
$ws.Range("G2:P2").Clear()

$testCases = @(
    "21.Create_Vendor_with_Questionnaire_with_Global_and_Local_and_Bank_with_Discard_JDE.xml",
    "20.Create_Vendor_with_Questionnaire_with_Global_and_Local_and_Bank_JDE.xml",
    "25Create_Vendor_with_Questionnaire_with_Global_and_Local_and_Bank_NAV.xml",
    "97.Create_Vendor_with_Questionnaire_banklocalonly_Global2.xml",
    "4.Change_vendor_global_&_local_&_Bank.xml"
)

$row = 2
foreach ($tc in $testCases) {
    $ws.Cells.Item($row, 1).Value = [string]($row - 1)
    $ws.Cells.Item($row, 2).Value = $tc
    $ws.Cells.Item($row, 3).Value = "PA, Panama"
    $ws.Cells.Item($row, 4).Value = "002"
    $ws.Cells.Item($row, 5).Value = "PAB, Panamanian Balboa"
    $ws.Cells.Item($row, 6).Value = "PAB1"
    $row++
}

$wb.Save()
